# Generate Report for Handback
# - Status changes from "Ready for handoff" to "Handback transform failed"
#   for the efa61409-... file (row 3) on every sheet that shows it
#   (Overview, zh-cn, de-de).
# - A new "Error Detail" (column K) message is recorded for row 3 on the
#   zh-cn and de-de sheets, describing the handback/handoff filename
#   mismatch that caused the transform to fail.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# Overview sheet: the efa61409 row's zh-cn (B) and de-de (C) status
# columns both carried the old status text.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# zh-cn sheet: Status column (C) + new Error Detail (K)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = "Handback file name: seqpqpir.zon is different with handoff file name: efa61409-2455-4ba1-bd46-6dc10031e533.b109b1e87d13141b0ba54b057b2994f8ce2b74b3.zh-cn."

# de-de sheet: Status column (C) + new Error Detail (K)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = "Handback file name: seqpqpir.zon is different with handoff file name: efa61409-2455-4ba1-bd46-6dc10031e533.b109b1e87d13141b0ba54b057b2994f8ce2b74b3.de-de."
